$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.272.93"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.14%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.676.87"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.78%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.007"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.22%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.77%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5287"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.97%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.007"
$ws.Range("D7").Style = "Normal"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2689"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.18%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06469"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.51%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.92"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.60%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07513"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.36%  "

$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.698.53"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.04%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.514"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.58%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5778"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.48%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000008494"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.62%  "

$ws.Range("E16").Value = "  +1.03%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.316.59"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.10%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.916"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.22%  "

$ws.Range("E19").Value = "  +0.16%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.85"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.75%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "189.98"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.66%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.195"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.20%  "

$ws.Range("E23").Value = "  +0.17%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "144.73"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.05%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1276"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +7.54%  "

$ws.Range("E26").Value = "  +2.96%  "

$ws.Range("E27").Value = "  +1.19%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06476"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.51%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.363"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.91%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.318"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.36%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.585"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.83%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.588"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.57%  "

$ws.Range("E33").Value = "  +1.67%  "

$ws.Range("E34").Value = "  +1.55%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6189"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.38%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.405"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.51%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.733"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.25%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.280"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.46%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.116.21"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.09%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01622"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.16%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8730"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.72%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.014"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.64%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.49"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.28%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.827.01"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.87%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000112"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.67%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.91"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.55%  "

$ws.Range("E47").Value = "  +0.21%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.169"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.12%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05263"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.13%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4291"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.02%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.048"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.94%  "
